# Commit: "Add registration info * Add symposium and young stats session"
#
# 1) A new "Satellite Symposium" sheet is added as the last sheet. It takes
#    over the "STRATOS" row that used to live at the bottom of the
#    "Featured Sessions" sheet (row 6), plus an extra (empty, wrap-formatted)
#    spacer row -- this reads like registration/info text was pasted below
#    the table and then cleared, leaving just the formatting behind.
# 2) The vacated row 6 on "Featured Sessions" is reused for the new
#    "Young Statisticians Sessions and Panel Discussion" entry.
# 3) The new sheet becomes the active / selected tab.

$wb = $excel.ActiveWorkbook

$featured = $wb.Worksheets.Item("Featured Sessions")

# Keep formatting references before we start moving things around.
$headerFormat = $featured.Range("A1:B1")
$wrapFormat   = $featured.Range("A5")

# --- Featured Sessions: row 6 becomes the Young Statisticians entry -------
$featured.Range("A6").Value = "Young Statisticians Sessions and Panel Discussion"
$featured.Range("B6").Value = "Andrea Berghold, Stefanie Peschel"

# --- New "Satellite Symposium" sheet, inserted after the last sheet ------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$satellite = $wb.Worksheets.Add($null, $lastSheet)
$satellite.Name = "Satellite Symposium"

$satellite.Range("A1").Value = "Title"
$satellite.Range("B1").Value = "Organizer"

$satellite.Range("A2").Value = "Ten years of the STRengthening Analytical Thinking for Observational Studies (STRATOS) initiative – progress and looking to the future"
$satellite.Range("B2").Value = "Ruth Keogh, Willi Sauerbrei"

# Match the header + wrap-text formatting used on the other sheets.
$headerFormat.Copy()
$satellite.Range("A1:B1").PasteSpecial(-4122)

$wrapFormat.Copy()
$satellite.Range("A5").PasteSpecial(-4122)

$satellite.Columns.Item(1).ColumnWidth = 108.69
$satellite.Columns.Item(2).ColumnWidth = 25.38

$satellite.PageSetup.Orientation = 1

$excel.CutCopyMode = $false

# --- Selections / active tab ---------------------------------------------
$featured.Activate()
$featured.Range("A7:B22").Select()

$satellite.Activate()
$satellite.Range("A7").Select()

Write-Host "done"
